# Updates cryptos list values (Price and Volume(1h) columns) to match
# the latest scrape, as produced by the GitHub Actions scheduled job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.689.45"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "2.304.01"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'302.34"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "'98.99"
$ws.Range("E6").Value = "  -5.85%  "
$ws.Range("D7").Value = "'0.501"
$ws.Range("E7").Value = "  -4.85%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").Value = "'34.65"
$ws.Range("E10").Value = "  -3.40%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'6.70"
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").Value = "2.657.65"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'15.61"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "2.298.25"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'0.796"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "42.593.29"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").Value = "'11.59"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "'6.04"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").Value = "'67.78"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'235.17"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").Value = "'1.95"
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").Value = "'2.51"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'24.65"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Value = "'34.29"
$ws.Range("E29").Value = "  -5.16%  "
$ws.Range("D30").Value = "'163.85"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "'9.11"
$ws.Range("E31").Value = "  -4.94%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'5.00"
$ws.Range("E33").Value = "  -4.73%  "
$ws.Range("E34").Value = "  -4.62%  "
$ws.Range("D35").Value = "'4.46"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'16.69"
$ws.Range("E36").Value = "  -9.06%  "
$ws.Range("D37").Value = "'2.88"
$ws.Range("E37").Value = "  -3.94%  "
$ws.Range("D38").Value = "'0.0694"
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").Value = "'0.0999"
$ws.Range("E40").Value = "  -5.63%  "
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").Value = "'2.48"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "1.962.45"
$ws.Range("D44").Value = "'0.0279"
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("D45").Value = "'18.42"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").Value = "'10.17"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").Value = "'2.87"
$ws.Range("E47").Value = "  -6.63%  "
$ws.Range("D48").Value = "'55.27"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("D49").Value = "2.527.10"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("D51").Value = "'4.64"
$ws.Range("E51").Value = "  -0.09%  "
